$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("I2").Value = 0.05501054111033029
$ws.Range("K2").Value = 169.56

# Row 3 updates
$ws.Range("H3").Value = 0.7897667942916812
$ws.Range("I3").Value = 0.05988760098349139
$ws.Range("K3").Value = 184
$ws.Range("R3").Value = 31
$ws.Range("T3").Value = 132
$ws.Range("U3").Value = 234
$ws.Range("W3").Value = 2816
$ws.Range("Y3").Value = 2715
$ws.Range("Z3").Value = 2613
$ws.Range("AG3").Value = 0.989111
$ws.Range("AI3").Value = 0.953635
$ws.Range("AJ3").Value = 0.917808
